# "1st changes of mifos to finflux"
#
# On the "Repayment schedule" sheet a new (blank) column is inserted
# immediately before column N, pushing the old N/O/P columns (Late,
# heading, Outstanding) one slot to the right (-> O/P/Q). The newly
# inserted column inherits the column width that used to belong to the
# column immediately to its left (M), exactly like Excel's native
# "Insert Column" command does.
#
# The previously-active sheet ("Transactions") is deselected and
# "Repayment schedule" becomes the active sheet/tab, with cell R13
# selected on it (and the 100% zoom level carried over).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N,O,P -> O,P,Q).
$ws.Columns("N").Insert()

# Match Excel's own behaviour of carrying the left-neighbour's column
# width onto the freshly inserted column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab (this also clears the
# previous tabSelected/active-cell state that used to live on the
# "Transactions" sheet) and select R13 on it.
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 100
$ws.Range("R13").Select() | Out-Null
